$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New URLs that need to be appended before the last existing row (which will shift down)
$newUrls = @(
    "https://stackoverflow.co/teams/",
    "https://stackoverflow.co/talent/",
    "https://stackoverflow.co/advertising/",
    "https://stackoverflow.co/labs/",
    "https://stackoverflow.co/"
)

# Insert new rows before row 6, pushing the existing row 6 (and below) down
$insertRange = $ws.Range("A6:A10")
$insertRange.Insert()

# Fill the newly inserted rows with the new URL values
for ($i = 0; $i -lt $newUrls.Length; $i++) {
    $rowIndex = 6 + $i
    $ws.Cells.Item($rowIndex, 1).Value = $newUrls[$i]
}
